$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor data-refresh corrections to existing "Ventas SOFOFA (base 2014=100)" values (column B)
$ws.Range("B182").Value = 95.45999999999999
$ws.Range("B186").Value = 99.06999999999999
$ws.Range("B187").Value = 98.69
$ws.Range("B188").Value = 99.12
$ws.Range("B221").Value = 97.54000000000001
$ws.Range("B235").Value = 97.37
$ws.Range("B239").Value = 96.72
$ws.Range("B240").Value = 106.57
$ws.Range("B241").Value = 107.31
$ws.Range("B246").Value = 95.70999999999999
$ws.Range("B258").Value = 105.99
$ws.Range("B260").Value = 112.68
$ws.Range("B261").Value = 116.33

# Append new row with the latest data point (01-09-2021)
$ws.Range("A262").NumberFormat = "@"
$ws.Range("A262").Value = "01-09-2021"
$ws.Range("A262").ClearFormats()
$ws.Range("B262").Value = 102.39
